$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param([string]$CellRef, [string]$Text)
    $rng = $ws.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.ClearFormats()
}

Set-CellText 'D2' '60.853.86'
Set-CellText 'E2' '  -4.66%  '
Set-CellText 'D3' '2.458.13'
Set-CellText 'E3' '  -6.02%  '
Set-CellText 'E4' '  +0.04%  '
Set-CellText 'D5' '545.85'
Set-CellText 'E5' '  -5.02%  '
Set-CellText 'D6' '145.05'
Set-CellText 'E6' '  -7.32%  '
Set-CellText 'E8' '  -4.39%  '
Set-CellText 'D9' '2.454.76'
Set-CellText 'E9' '  -6.02%  '
Set-CellText 'D10' '0.107'
Set-CellText 'E10' '  -10.19%  '
Set-CellText 'E11' '  -2.01%  '
Set-CellText 'D12' '5.34'
Set-CellText 'E12' '  -8.58%  '
Set-CellText 'D13' '0.351'
Set-CellText 'E13' '  -8.13%  '
Set-CellText 'D14' '25.89'
Set-CellText 'E14' '  -8.31%  '
Set-CellText 'D15' '2.899.38'
Set-CellText 'E16' '  -10.41%  '
Set-CellText 'D17' '60.775.88'
Set-CellText 'E17' '  -4.54%  '
Set-CellText 'D18' '2.449.26'
Set-CellText 'E18' '  -5.26%  '
Set-CellText 'D19' '11.00'
Set-CellText 'E19' '  -8.59%  '
Set-CellText 'D20' '6.93'
Set-CellText 'E20' '  -8.49%  '
Set-CellText 'D21' '4.16'
Set-CellText 'E21' '  -8.32%  '
Set-CellText 'D22' '317.27'
Set-CellText 'E22' '  -7.73%  '
Set-CellText 'E23' '  -0.07%  '
Set-CellText 'D24' '62.98'
Set-CellText 'E24' '  -6.58%  '
Set-CellText 'E25' '  -5.87%  '
Set-CellText 'D26' '0.0₃0976'
Set-CellText 'E26' '  -10.45%  '
Set-CellText 'D27' '2.576.45'
Set-CellText 'E27' '  -5.65%  '
Set-CellText 'D28' '0.999'
Set-CellText 'E28' '  -0.11%  '
Set-CellText 'D29' '1.48'
Set-CellText 'E29' '  -5.96%  '
Set-CellText 'D30' '530.21'
Set-CellText 'E30' '  -11.23%  '
Set-CellText 'E31' '  -10.00%  '
Set-CellText 'D32' '7.64'
Set-CellText 'E32' '  -3.59%  '
Set-CellText 'D33' '0.149'
Set-CellText 'E33' '  -8.13%  '
Set-CellText 'E34' '  -8.82%  '
Set-CellText 'E35' '  -10.09%  '
Set-CellText 'E36' '  -11.88%  '
Set-CellText 'B37' 'NEARProtocol'
Set-CellText 'C37' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-CellText 'D37' '4.84'
Set-CellText 'E37' '  -10.57%  '
Set-CellText 'B38' 'FirstDigitalUSD'
Set-CellText 'C38' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-CellText 'D38' '0.998'
Set-CellText 'E38' '  -0.07%  '
Set-CellText 'E39' '  -6.98%  '
Set-CellText 'D40' '18.25'
Set-CellText 'E40' '  -7.70%  '
Set-CellText 'D41' '144.35'
Set-CellText 'E41' '  -6.58%  '
Set-CellText 'E42' '  -0.10%  '
Set-CellText 'E43' '  -9.90%  '
Set-CellText 'E44' '  -3.95%  '
Set-CellText 'D45' '2.28'
Set-CellText 'E45' '  -10.63%  '
Set-CellText 'D46' '146.27'
Set-CellText 'E46' '  -7.07%  '
Set-CellText 'E47' '  -9.08%  '
Set-CellText 'E48' '  -13.25%  '
Set-CellText 'D49' '0.0528'
Set-CellText 'E49' '  -10.66%  '
Set-CellText 'D50' '0.581'
Set-CellText 'E50' '  -7.84%  '
Set-CellText 'D51' '0.0938'
Set-CellText 'E51' '  -6.37%  '
